# Weekly NYPD CompStat data refresh: bump the report week / volume number,
# and update this week's crime-statistics table (rows 16-30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text: "Volume 30   Number  45" -> "...  46"
# and the "Report Covering the Week 11/6/2023 Through 11/12/2023" line
# -> "...11/13/2023 ... 11/19/2023"
# These are rich-text shared strings; Excel re-assembles them as a single
# concatenated string when you set .Value on the cell (runs collapse, but
# the run formatting in the source file was identical for every run so
# nothing visually changes).
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  46"
$ws.Range("C9").Value = "Report Covering the Week  11/13/2023  Through  11/19/2023"

# ---------------------------------------------------------------------
# Helper: PasteSpecial-formats (xlPasteFormats = -4122) lets us copy just
# the number format / style from a donor cell onto a target cell whose
# data type changed (number <-> text placeholder), without introducing a
# brand-new style entry.
# ---------------------------------------------------------------------

function Set-NumValue($addr, $val) {
    $ws.Range($addr).Value = $val
}

function Set-TextLikeValue($addr, $text, $donorAddr) {
    # Force text storage for a numeric-looking placeholder ("0", "***.*")
    # by entering it with a leading apostrophe, then copy the format from
    # a donor cell that already carries the correct text style.
    $ws.Range($addr).Value = "'" + $text
    $ws.Range($donorAddr).Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}

function Set-NumValueWithFormat($addr, $val, $donorAddr) {
    # Force numeric storage (coming from a text placeholder cell) and pick
    # up the correct numeric style from a donor cell.
    $ws.Range($addr).Value = $val
    $ws.Range($donorAddr).Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}

# ---------------------------------------------------------------------
# Row 16
# ---------------------------------------------------------------------
Set-NumValue "C16" 6
Set-NumValue "D16" 9
Set-NumValue "E16" -33.333333333333
Set-NumValue "F16" 24
Set-NumValue "G16" 23
Set-NumValue "H16" 4.347826086956
Set-NumValue "I16" 174
Set-NumValue "J16" 229
Set-NumValue "K16" -24.017467248908
Set-NumValue "L16" 4.819277108433
Set-NumValue "M16" 31.818181818181
Set-NumValue "N16" -77.373211963589

# ---------------------------------------------------------------------
# Row 17
# ---------------------------------------------------------------------
Set-NumValue "C17" 1
Set-NumValue "D17" 1
Set-NumValue "E17" 0
Set-NumValue "F17" 10
Set-NumValue "G17" 12
Set-NumValue "H17" -16.666666666666
Set-NumValue "I17" 177
Set-NumValue "J17" 166
Set-NumValue "K17" 6.626506024096
Set-NumValue "L17" 9.259259259259
Set-NumValue "M17" 108.235294117647
Set-NumValue "N17" -35.869565217391

# ---------------------------------------------------------------------
# Row 18
# ---------------------------------------------------------------------
Set-NumValue "C18" 12
Set-NumValue "D18" 5
Set-NumValue "E18" 140
Set-NumValue "F18" 31
Set-NumValue "G18" 34
Set-NumValue "H18" -8.823529411764
Set-NumValue "I18" 263
Set-NumValue "J18" 390
Set-NumValue "K18" -32.564102564102
Set-NumValue "L18" 22.325581395348
Set-NumValue "M18" 64.375
Set-NumValue "N18" -63.972602739726

# ---------------------------------------------------------------------
# Row 19
# ---------------------------------------------------------------------
Set-NumValue "C19" 20
Set-NumValue "D19" 31
Set-NumValue "E19" -35.483870967741
Set-NumValue "F19" 149
Set-NumValue "G19" 158
Set-NumValue "H19" -5.696202531645
Set-NumValue "I19" 1187
Set-NumValue "J19" 1233
Set-NumValue "K19" -3.730738037307
Set-NumValue "L19" 58.689839572192
Set-NumValue "M19" 27.909482758620
Set-NumValue "N19" -48.503253796095

# ---------------------------------------------------------------------
# Row 20 (C/D/E untouched by this week's refresh)
# ---------------------------------------------------------------------
Set-NumValue "F20" 2
Set-NumValue "G20" 4
Set-NumValue "H20" -50
Set-NumValue "I20" 43
Set-NumValue "J20" 55
Set-NumValue "K20" -21.818181818181
Set-NumValue "L20" -8.510638297872
Set-NumValue "M20" 4.878048780487
Set-NumValue "N20" -93.740902474526

# ---------------------------------------------------------------------
# Row 21 (TOTAL row)
# ---------------------------------------------------------------------
Set-NumValue "C21" 39
Set-NumValue "D21" 47
Set-NumValue "E21" -17.021276595744
Set-NumValue "F21" 217
Set-NumValue "G21" 232
Set-NumValue "H21" -6.465517241379
Set-NumValue "I21" 1853
Set-NumValue "J21" 2087
Set-NumValue "K21" -11.212266411116
Set-NumValue "L21" 37.564959168522
Set-NumValue "M21" 36.450662739322
Set-NumValue "N21" -61.226197949361

# ---------------------------------------------------------------------
# Row 22 (Transit): C/D/E flip from numeric 1/1/0 to text placeholders
# "0"/"0"/"***.*" - no crimes of this type this week.
# ---------------------------------------------------------------------
Set-TextLikeValue "C22" "0" "C15"
Set-TextLikeValue "D22" "0" "D15"
Set-TextLikeValue "E22" "***.*" "E15"

# ---------------------------------------------------------------------
# Row 24 (Petit Larceny)
# ---------------------------------------------------------------------
Set-NumValue "C24" 34
Set-NumValue "D24" 51
Set-NumValue "E24" -33.333333333333
Set-NumValue "F24" 164
Set-NumValue "G24" 152
Set-NumValue "H24" 7.894736842105
Set-NumValue "I24" 1869
Set-NumValue "J24" 1757
Set-NumValue "K24" 6.374501992031
Set-NumValue "L24" 61.398963730569
Set-NumValue "M24" 39.477611940298

# ---------------------------------------------------------------------
# Row 25 (Misd. Assault)
# ---------------------------------------------------------------------
Set-NumValue "C25" 7
Set-NumValue "D25" 12
Set-NumValue "E25" -41.666666666666
Set-NumValue "F25" 34
Set-NumValue "G25" 47
Set-NumValue "H25" -27.659574468085
Set-NumValue "I25" 383
Set-NumValue "J25" 380
Set-NumValue "K25" 0.789473684210
Set-NumValue "L25" 35.335689045936
Set-NumValue "M25" 64.377682403433

# ---------------------------------------------------------------------
# Row 27 (Other Sex Crimes): C flips numeric->text "0"; D/E flip the
# other way, text placeholders -> real numbers, now that data exists.
# Grab the numeric donor styles from this same row (C27 is s15-style
# before we touch it, H27 is s16-style) before we overwrite C27.
# ---------------------------------------------------------------------
Set-NumValueWithFormat "D27" 3 "C27"
Set-NumValueWithFormat "E27" -100 "H27"
Set-TextLikeValue "C27" "0" "C15"
Set-NumValue "F27" 4
Set-NumValue "G27" 6
Set-NumValue "H27" -33.333333333333
Set-NumValue "J27" 70
Set-NumValue "K27" -17.142857142857

# ---------------------------------------------------------------------
# Row 30 (Hate Crimes): D/E flip from numeric to text placeholders
# "0"/"***.*" ; I/K/L get refreshed numbers.
# ---------------------------------------------------------------------
Set-TextLikeValue "D30" "0" "C15"
Set-TextLikeValue "E30" "***.*" "E15"
Set-NumValue "I30" 9
Set-NumValue "K30" -25
Set-NumValue "L30" -35.714285714285
